$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of kaspa buy data appended after the 2025-05-20 run.
# The Date column in this workbook's recent rows is stored as plain text
# (e.g. "05/17/2025"), not an Excel date serial, so force text entry with a
# leading apostrophe and reset the cell style to Normal to avoid Excel's
# automatic date-number reformatting.
$ws.Range("A26").Value = "'05/20/2025"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").Value = 459.3410000000003
$ws.Range("C26").Value = 0.1088515939138896
$ws.Range("D26").Value = 50
